$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.391.36'
$ws.Range('E2').Value = '  -3.27%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.667.73'
$ws.Range('E3').Value = '  -2.31%  '
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '219.05'
$ws.Range('E5').Value = '  -1.95%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5180'
$ws.Range('E6').Value = '  -2.50%  '
$ws.Range('E7').Value = '  +0.50%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.06464'
$ws.Range('E8').Value = '  -1.79%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2579'
$ws.Range('E9').Value = '  -2.88%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.99'
$ws.Range('E10').Value = '  -3.93%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07671'
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.684.85'
$ws.Range('E12').Value = '  -1.29%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.346'
$ws.Range('E13').Value = '  -4.72%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.895.11'
$ws.Range('E14').Value = '  -2.52%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5561'
$ws.Range('E15').Value = '  -2.78%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0₅8056'
$ws.Range('E16').Value = '  -1.29%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '64.75'
$ws.Range('E17').Value = '  -4.37%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '26.438.66'
$ws.Range('E18').Value = '  -3.31%  '
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '210.13'
$ws.Range('E20').Value = '  -2.61%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.425'
$ws.Range('E21').Value = '  -5.10%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.13'
$ws.Range('E22').Value = '  -2.83%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.896'
$ws.Range('E23').Value = '  -1.17%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.008'
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '145.09'
$ws.Range('E25').Value = '  +2.49%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.728'
$ws.Range('E26').Value = '  -2.17%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1169'
$ws.Range('E27').Value = '  -3.73%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.010'
$ws.Range('E28').Value = '  -3.42%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.82'
$ws.Range('E29').Value = '  -2.79%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05254'
$ws.Range('E30').Value = '  -2.87%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.263'
$ws.Range('E31').Value = '  -2.16%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.383'
$ws.Range('E32').Value = '  -3.36%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.223'
$ws.Range('E33').Value = '  -5.89%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.578'
$ws.Range('E34').Value = '  -4.04%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.760'
$ws.Range('E35').Value = '  -3.93%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.377'
$ws.Range('E36').Value = '  -1.17%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9299'
$ws.Range('E37').Value = '  -1.85%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5724'
$ws.Range('E38').Value = '  -2.13%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.152.25'
$ws.Range('E39').Value = '  +10.49%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01607'
$ws.Range('E40').Value = '  -1.28%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.007'
$ws.Range('E41').Value = '  +0.54%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.8456'
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.655'
$ws.Range('E43').Value = '  -3.44%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '100.32'
$ws.Range('E44').Value = '  -0.52%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.805.29'
$ws.Range('E45').Value = '  -2.49%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0₈110'
$ws.Range('E46').Value = '  -6.62%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4496'
$ws.Range('E47').Value = '  -0.05%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '56.06'
$ws.Range('E48').Value = '  -3.15%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.007'
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.946'
$ws.Range('E50').Value = '  -1.98%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05115'
$ws.Range('E51').Value = '  -2.44%  '
